# Adds a new "Swiss" worksheet (test data for the Switzerland market),
# cloned from the existing "Czech" sheet, with the market name / ticket
# reference updated, and makes it the active sheet/tab.

$wb = $excel.ActiveWorkbook

# Clone the "Czech" sheet (same layout/styles/merged cells) and place the
# copy right after it - this becomes the 4th sheet, matching the diff's
# new <sheet name="Swiss" .../> entry appended at the end of the list.
$czech = $wb.Worksheets.Item("Czech")
[void]$czech.Copy([System.Reflection.Missing]::Value, $czech)

# The newly inserted sheet is now the last sheet in the workbook.
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Update the market name and the associated ticket/reference number.
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2642"

# Match the saved selection/active-cell state from the diff:
#  - new "Swiss" sheet ends up as the active tab, with B14 selected
#  - "Czech" sheet is left with the whole sheet selected (no active tab)
[void]$swiss.Range("B14").Select()
[void]$czech.Select()
[void]$czech.Cells.Select()
[void]$swiss.Activate()
